# Auto-generated edit script applying the Ifrit_Profits.xlsx commit diff
# to the corresponding profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1109.7858
$ws.Range("I19").Value = 1225
$ws.Range("J19").Value = 1023.375
$ws.Range("K19").Value = 1225
$ws.Range("L19").Value = 1023.375
$ws.Range("M19").Value = -1050
$ws.Range("N19").Value = -1373.375

$ws.Range("H40").Value = 1075.7142
$ws.Range("I40").Value = 1010
$ws.Range("J40").Value = 1240
$ws.Range("K40").Value = 1010
$ws.Range("L40").Value = 1240
$ws.Range("M40").Value = -835
$ws.Range("N40").Value = -1590

$ws.Range("H64").Value = 94145.45
$ws.Range("I64").Value = 4200
$ws.Range("J64").Value = 127875
$ws.Range("K64").Value = 4200
$ws.Range("L64").Value = 127875
$ws.Range("M64").Value = -3952
$ws.Range("N64").Value = -128371

$ws.Range("H67").Value = 94145.45
$ws.Range("I67").Value = 4200
$ws.Range("J67").Value = 127875
$ws.Range("K67").Value = 4200
$ws.Range("L67").Value = 127875
$ws.Range("M67").Value = -3342
$ws.Range("N67").Value = -129591

$ws.Range("H112").Value = 50001308
$ws.Range("I112").Value = 506
$ws.Range("J112").Value = 66668240
$ws.Range("K112").Value = 1518
$ws.Range("L112").Value = 200004720
$ws.Range("M112").Value = -410
$ws.Range("N112").Value = -200006936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2675914
$ws.Range("I2").Value = 2573.875
$ws.Range("J2").Value = 9804821
$ws.Range("K2").Value = 2573.875
$ws.Range("L2").Value = 9804821
$ws.Range("M2").Value = -2460.875
$ws.Range("N2").Value = -9805047

$ws.Range("H32").Value = 4313.276
$ws.Range("I32").Value = 4745.4146
$ws.Range("J32").Value = 3271.0588
$ws.Range("K32").Value = 4745.4146
$ws.Range("L32").Value = 3271.0588
$ws.Range("M32").Value = -4458.4146

$ws.Range("H63").Value = 4450.6665
$ws.Range("I63").Value = 2011.4
$ws.Range("J63").Value = 7499.75
$ws.Range("K63").Value = 2011.4
$ws.Range("L63").Value = 7499.75
$ws.Range("M63").Value = -1325.4

$ws.Range("H66").Value = 4450.6665
$ws.Range("I66").Value = 2011.4
$ws.Range("J66").Value = 7499.75
$ws.Range("K66").Value = 10057
$ws.Range("L66").Value = 37498.75
$ws.Range("M66").Value = -6625

$ws.Range("H97").Value = 542.1111
$ws.Range("I97").Value = 542.1111
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 542.1111
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -46.11109999999996
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 1463.7646
$ws.Range("I102").Value = 1478.125
$ws.Range("J102").Value = 1234
$ws.Range("K102").Value = 1478.125
$ws.Range("L102").Value = 1234
$ws.Range("M102").Value = 143.875
$ws.Range("N102").Value = -4478

$ws.Range("H116").Value = 2675914
$ws.Range("I116").Value = 2573.875
$ws.Range("J116").Value = 9804821
$ws.Range("K116").Value = 2573.875
$ws.Range("L116").Value = 9804821
$ws.Range("M116").Value = -279.875
$ws.Range("N116").Value = -9809409

$ws.Range("H122").Value = 2167.375
$ws.Range("I122").Value = 2120.5
$ws.Range("J122").Value = 2495.5
$ws.Range("K122").Value = 6361.5
$ws.Range("L122").Value = 7486.5
$ws.Range("M122").Value = -3911.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2675914
$ws.Range("I3").Value = 2573.875
$ws.Range("J3").Value = 9804821
$ws.Range("K3").Value = 2573.875
$ws.Range("L3").Value = 9804821
$ws.Range("M3").Value = -2459.875
$ws.Range("N3").Value = -9805049

$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2549
$ws.Range("N94").ClearContents()

$ws.Range("H99").Value = 946.41174
$ws.Range("I99").Value = 1008.9091
$ws.Range("J99").Value = 831.8333
$ws.Range("K99").Value = 1008.9091
$ws.Range("L99").Value = 831.8333
$ws.Range("M99").Value = 489.0909
$ws.Range("N99").Value = -3827.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3250
$ws.Range("I16").Value = 3250
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3250
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2963

$ws.Range("H62").Value = 3158.4666
$ws.Range("I62").Value = 3037.8
$ws.Range("J62").Value = 3399.8
$ws.Range("K62").Value = 3037.8
$ws.Range("L62").Value = 3399.8
$ws.Range("M62").Value = -2413.8
$ws.Range("N62").Value = -4647.8

$ws.Range("H65").Value = 3158.4666
$ws.Range("I65").Value = 3037.8
$ws.Range("J65").Value = 3399.8
$ws.Range("K65").Value = 15189
$ws.Range("L65").Value = 16999
$ws.Range("M65").Value = -12069
$ws.Range("N65").Value = -23239

$ws.Range("H113").Value = 3250
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1000
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2441
$ws.Range("N57").ClearContents()

$ws.Range("H63").Value = 2666.6667
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2666.6667
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 8000.000100000001
$ws.Range("N63").Value = -9498.000100000001
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 2666.6667
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 2666.6667
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 24000.0003
$ws.Range("N66").Value = -31488.0003
$ws.Range("M66").ClearContents()

$ws.Range("H94").Value = 2705.8235
$ws.Range("I94").Value = 1999.6666
$ws.Range("J94").Value = 2857.1428
$ws.Range("K94").Value = 5998.9998
$ws.Range("L94").Value = 8571.428400000001
$ws.Range("M94").Value = -5322.9998
$ws.Range("N94").Value = -9923.428400000001

$ws.Range("H95").Value = 3150
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 3150
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 9450
$ws.Range("N95").Value = -13568

$ws.Range("H101").Value = 7900
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7900
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23700
$ws.Range("N101").Value = -28568

$ws.Range("H110").Value = 1933.3334
$ws.Range("I110").Value = 1933.3334
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5800.0002
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1710.0002

$ws.Range("H131").Value = 1820947.6
$ws.Range("I131").Value = 4052
$ws.Range("J131").Value = 2566340.8
$ws.Range("K131").Value = 12156
$ws.Range("L131").Value = 7699022.399999999
$ws.Range("M131").Value = -7116
$ws.Range("N131").Value = -7709102.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6000.8887
$ws.Range("I70").Value = 6286.857
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6286.857
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -6016.857
$ws.Range("N70").Value = -5540

$ws.Range("H73").Value = 6000.8887
$ws.Range("I73").Value = 6286.857
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6286.857
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -5350.857
$ws.Range("N73").Value = -6872

$ws.Range("H80").Value = 121600.4
$ws.Range("I80").Value = 2499.3333
$ws.Range("J80").Value = 172643.72
$ws.Range("K80").Value = 2499.3333
$ws.Range("L80").Value = 172643.72
$ws.Range("M80").Value = -1501.3333
$ws.Range("N80").Value = -174639.72

$ws.Range("H83").Value = 121600.4
$ws.Range("I83").Value = 2499.3333
$ws.Range("J83").Value = 172643.72
$ws.Range("K83").Value = 12496.6665
$ws.Range("L83").Value = 863218.6
$ws.Range("M83").Value = -7504.666499999999
$ws.Range("N83").Value = -873202.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H127").Value = 46021
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 46021
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 46021
$ws.Range("N127").Value = -55941
